# Weekly data refresh: insert one new observation row for Espinaca
# (Mercado Mayorista Lo Valledor de Santiago) ahead of the existing
# row 606, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 606 (pushes old rows 606..682 down to 607..683)
$ws.Rows.Item(606).Insert()

# Populate the new row 606 with the latest weekly observation
$ws.Cells.Item(606, 1).Value  = 6
$ws.Cells.Item(606, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(606, 3).Value  = "Metropolitana"
$ws.Cells.Item(606, 4).Value  = 44946
$ws.Cells.Item(606, 5).Value  = 13
$ws.Cells.Item(606, 6).Value  = 100112012
$ws.Cells.Item(606, 7).Value  = "Espinaca"
$ws.Cells.Item(606, 8).Value  = "Sin especificar"
$ws.Cells.Item(606, 9).Value  = "Primera"
$ws.Cells.Item(606, 10).Value = 610
$ws.Cells.Item(606, 11).Value = 5000
$ws.Cells.Item(606, 12).Value = 6000
$ws.Cells.Item(606, 13).Value = 5410
$ws.Cells.Item(606, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(606, 15).Value = "Región Metropolitana"
$ws.Cells.Item(606, 16).Value = 541
$ws.Cells.Item(606, 17).Value = 10
$ws.Cells.Item(606, 18).Value = "Hortaliza"
